# Edit: (1) re-style the table on slide 5 with the built-in
# "Themed Style 1 - Accent 6" table style, and (2) switch the deck's
# design theme over to the default Office theme palette (the deck
# shipped with both an "Integral" theme, used by the slide master, and
# an "Office Theme" palette, used only by the notes master -- this
# commit makes the Office palette the active one).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Table style id -> built-in "Themed Style 1 - Accent 6"
# ---------------------------------------------------------------
$s = $p.Slides.Item(5)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{EAFE2B55-1022-42B4-9315-E4FF975B93B3}")
    }
}

# ---------------------------------------------------------------
# 2) Re-colour the active theme to the standard Office palette
# ---------------------------------------------------------------
function Set-RGB($themeColors, $index, $r, $g, $b) {
    $themeColors.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

Set-RGB $themeColors 1  0x00 0x00 0x00   # dk1
Set-RGB $themeColors 2  0xFF 0xFF 0xFF   # lt1
Set-RGB $themeColors 3  0x44 0x54 0x6A   # dk2
Set-RGB $themeColors 4  0xE7 0xE6 0xE6   # lt2
Set-RGB $themeColors 5  0x5B 0x9B 0xD5   # accent1
Set-RGB $themeColors 6  0xED 0x7D 0x31   # accent2
Set-RGB $themeColors 7  0xA5 0xA5 0xA5   # accent3
Set-RGB $themeColors 8  0xFF 0xC0 0x00   # accent4
Set-RGB $themeColors 9  0x44 0x72 0xC4   # accent5
Set-RGB $themeColors 10 0x70 0xAD 0x47   # accent6
Set-RGB $themeColors 11 0x05 0x63 0xC1   # hlink
Set-RGB $themeColors 12 0x95 0x4F 0x72   # folHlink
